$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New block of data below the existing "sum 5-2 sort" table (rows 24-33) ---
# Entry order matters: it controls the order new strings are appended to the
# shared-strings table, so cells are written in the same sequence the
# original author must have used (birth row first, then survival row,
# then the lower/upper labels, the header, and finally the summary rows).

# Row 28 : "birth" lower row (label only first, to seed the "birth" string)
$ws.Range("B28").Value = "birth"

# Row 26 : "survival" lower row (seeds the "survival" string)
$ws.Range("B26").Value = "survival"

# Row 26 continued : "lower" label + counts (seeds the "lower" string)
$ws.Range("C26").Value = "lower"
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 1

# Row 27 : survival "upper" row (seeds the "upper" string)
$ws.Range("C27").Value = "upper"
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 1

# Row 28 continued : birth "lower" counts
$ws.Range("C28").Value = "lower"
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 1

# Row 29 : birth "upper" row
$ws.Range("C29").Value = "upper"
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 1
$ws.Range("K29").Value = 1
$ws.Range("L29").Value = 1

# Row 24 : section header
$ws.Range("A24").Value = "values inclusive"

# Row 25 : copy of the H column totals (dist sums) from the table above
$ws.Range("D25").Formula = "=H5"
$ws.Range("E25").Formula = "=H6"
$ws.Range("F25").Formula = "=H7"
$ws.Range("G25").Formula = "=H8"
$ws.Range("H25").Formula = "=H9"
$ws.Range("I25").Formula = "=H10"
$ws.Range("J25").Formula = "=H11"
$ws.Range("K25").Formula = "=H12"
$ws.Range("L25").Formula = "=H13"

# Row 30 : survival configs summary
$ws.Range("B30").Value = "survival configs"
$ws.Range("D30").Value = 7
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 4
$ws.Range("I30").Value = 3
$ws.Range("M30").Formula = "=SUM(D30:L30)"

# Row 31 : birth configs summary
$ws.Range("B31").Value = "birth configs"
$ws.Range("E31").Value = 7
$ws.Range("F31").Value = 7
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = 5
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 3
$ws.Range("M31").Formula = "=SUM(D31:L31)"

# Row 33 : closing note
$ws.Range("M33").Value = "1024 possibilities"

# Restore the active-cell selection to match the author's final cursor position
$null = $ws.Range("I39").Select()
